# Uploaded Course Materials - Day 10.
# Remove the empty "Title 1" placeholder shape from slide 4
# (Title 1 / id=2, type="title") while leaving the slide's other
# shapes (Slide Number Placeholder 4, Content Placeholder 5) intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Shapes.Item("Title 1").Delete()
